$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into the Price column (D) while preserving the
# column's original text storage. Every Price cell in this sheet holds a
# literal text string (e.g. "57.20", "2.50", "44.049.46") rather than a
# number; a plain Range.Value assignment would let Excel auto-coerce a
# numeric-looking string into an actual number and silently drop a
# meaningful trailing zero (e.g. "57.20" -> 57.2) or mis-parse a multi-dot
# string. Forcing a Text number format for the write, then restoring the
# default "Normal" style, keeps the stored value exact text while leaving
# cell formatting identical to the original workbook.
function Set-PriceCell($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-PriceCell 'D2' '43.967.47'
$ws.Range('E2').Value = '  +0.26%  '
Set-PriceCell 'D3' '2.359.75'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('E5').Value = '  +1.21%  '
Set-PriceCell 'D6' '239.61'
$ws.Range('E6').Value = '  +1.07%  '
Set-PriceCell 'D7' '74.52'
$ws.Range('E8').Value = '  +0.01%  '
Set-PriceCell 'D9' '0.594'
$ws.Range('E9').Value = '  +10.16%  '
$ws.Range('E10').Value = '  +1.55%  '
Set-PriceCell 'D11' '57.20'
$ws.Range('E11').Value = '  +0.02%  '
Set-PriceCell 'D12' '32.13'
$ws.Range('E12').Value = '  +12.85%  '
Set-PriceCell 'D13' '7.29'
$ws.Range('E13').Value = '  +9.86%  '
Set-PriceCell 'D14' '0.108'
$ws.Range('E14').Value = '  +0.78%  '
Set-PriceCell 'D15' '2.708.76'
$ws.Range('E15').Value = '  +0.29%  '
Set-PriceCell 'D16' '16.66'
$ws.Range('E16').Value = '  +0.00%  '
Set-PriceCell 'D17' '0.904'
$ws.Range('E17').Value = '  +1.10%  '
Set-PriceCell 'D18' '2.360.51'
$ws.Range('E18').Value = '  +0.04%  '
Set-PriceCell 'D19' '43.888.53'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('E20').Value = '  +1.11%  '
Set-PriceCell 'D21' '6.77'
$ws.Range('E21').Value = '  +5.37%  '
Set-PriceCell 'D22' '77.03'
$ws.Range('E22').Value = '  -1.07%  '
Set-PriceCell 'D23' '256.68'
$ws.Range('E23').Value = '  +1.15%  '
Set-PriceCell 'D24' '1.97'
$ws.Range('E24').Value = '  +24.99%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('E26').Value = '  -1.94%  '
$ws.Range('E27').Value = '  -0.48%  '
Set-PriceCell 'D28' '10.72'
$ws.Range('E28').Value = '  +1.77%  '
Set-PriceCell 'D29' '2.26'
$ws.Range('E29').Value = '  -0.80%  '
Set-PriceCell 'D30' '22.76'
$ws.Range('E30').Value = '  +1.73%  '
Set-PriceCell 'D31' '175.08'
$ws.Range('E31').Value = '  +1.42%  '
$ws.Range('E32').Value = '  -2.54%  '
Set-PriceCell 'D33' '0.137'
$ws.Range('E33').Value = '  +4.06%  '
Set-PriceCell 'D34' '0.0758'
$ws.Range('E34').Value = '  +6.46%  '
$ws.Range('E35').Value = '  +1.47%  '
Set-PriceCell 'D36' '5.42'
$ws.Range('E36').Value = '  +5.13%  '
$ws.Range('E37').Value = '  -7.29%  '
$ws.Range('E38').Value = '  -2.65%  '
Set-PriceCell 'D39' '6.31'
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('E40').Value = '  +4.28%  '
$ws.Range('E41').Value = '  +16.82%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-PriceCell 'D42' '0.205'
$ws.Range('E42').Value = '  +13.42%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-PriceCell 'D43' '9.16'
$ws.Range('E43').Value = '  +3.82%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-PriceCell 'D44' '19.21'
$ws.Range('E44').Value = '  -2.25%  '
$ws.Range('E45').Value = '  +0.03%  '
Set-PriceCell 'D46' '4.73'
$ws.Range('E46').Value = '  +6.58%  '
Set-PriceCell 'D47' '58.10'
$ws.Range('E47').Value = '  +10.49%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-PriceCell 'D48' '2.50'
$ws.Range('E48').Value = '  +7.95%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-PriceCell 'D49' '1.24'
$ws.Range('E49').Value = '  +0.45%  '
$ws.Range('E50').Value = '  +1.23%  '
Set-PriceCell 'D51' '100.36'
$ws.Range('E51').Value = '  +2.72%  '
